$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.510.36"
$ws.Range("E2").Value = "  +4.90%  "
$ws.Range("D3").Value = "'1.602.75"
$ws.Range("E3").Value = "  +2.61%  "
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").Value = "'215.28"
$ws.Range("E5").Value = "  +2.30%  "
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("D8").Value = "'24.13"
$ws.Range("E8").Value = "  +9.61%  "
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("D11").Value = "'0.0891"
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("D12").Value = "'1.829.41"
$ws.Range("E12").Value = "  +2.50%  "
$ws.Range("D13").Value = "'1.592.27"
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("E15").Value = "  +3.58%  "
$ws.Range("D16").Value = "'28.523.82"
$ws.Range("E16").Value = "  +4.92%  "
$ws.Range("D17").Value = "'63.48"
$ws.Range("E17").Value = "  +2.63%  "
$ws.Range("D18").Value = "'232.78"
$ws.Range("E18").Value = "  +7.44%  "
$ws.Range("D19").Value = "'7.57"
$ws.Range("E19").Value = "  +1.99%  "
$ws.Range("D20").Value = "'0.0₃0713"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").Value = "'4.14"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("E23").Value = "  +2.58%  "
$ws.Range("D24").Value = "'1.97"
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("D25").Value = "'152.82"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").Value = "'15.34"
$ws.Range("E26").Value = "  +2.37%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("D32").Value = "'3.26"
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("D34").Value = "'1.423.73"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -3.81%  "
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("D38").Value = "'0.0167"
$ws.Range("E38").Value = "  +0.98%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.52"
$ws.Range("E39").Value = "  +8.22%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "'0.545"
$ws.Range("E40").Value = "  +2.65%  "
$ws.Range("D41").Value = "'0.825"
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("D42").Value = "'5.76"
$ws.Range("E42").Value = "  -2.67%  "
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").Value = "'1.85"
$ws.Range("E44").Value = "  +7.03%  "
$ws.Range("D45").Value = "'0.977"
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").Value = "'65.17"
$ws.Range("E46").Value = "  +1.30%  "
$ws.Range("D47").Value = "'1.740.72"
$ws.Range("E47").Value = "  +2.53%  "
$ws.Range("D48").Value = "'87.64"
$ws.Range("E48").Value = "  +2.43%  "
$ws.Range("D49").Value = "'2.14"
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("D50").Value = "'0.0₆0108"
$ws.Range("E50").Value = "  +6.22%  "
$ws.Range("D51").Value = "'0.0528"
$ws.Range("E51").Value = "  +0.52%  "
